$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: add "Popravni" in column H for Mervan Bronja (failing grade)
$ws.Range("H9").Value = "Popravni"

# Row 14: fill in scores for student 11
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 20
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = 20

# Row 15: fill in scores for student 12
$ws.Range("C15").Value = 18
$ws.Range("D15").Value = 20
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 5

# Update selection to match target state
$ws.Range("H13").Select()
